{"js": "// The second paragraph of the document (the byline with the author's\n// name, course code and ID number, separated by tabs) is center-aligned\n// and stripped of all of its run content (\"Sami Mansoor Alavi\", the\n// tabs, \"BESE 8B\", the tabs, and \"209433\"), while the existing\n// \"_GoBack\" bookmark that sits in the middle of that text is preserved\n// in place.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph index 0 is the \"Secure Wallet\" title; index 1 is the byline\n// paragraph that needs to be cleared and centered.\nconst byline = paragraphs.items[1];\n\n// Center the paragraph (adds <w:jc w:val=\"center\"/> to its pPr).\nbyline.alignment = Word.Alignment.centered;\n\n// Locate the (hidden) \"_GoBack\" bookmark that lives between the tabs\n// after the name and before \"BESE 8B\" so we can delete the text around\n// it without disturbing it.\nconst bookmarkRange = body.getBookmarkRange(\"_GoBack\");\nconst paragraphStart = byline.getRange(\"Start\");\nconst paragraphEnd = byline.getRange(\"End\");\n\nconst beforeBookmark = paragraphStart.expandTo(bookmarkRange);\nconst afterBookmark = bookmarkRange.expandTo(paragraphEnd);\n\n// Clear the text on both sides of the bookmark (order matters: clear\n// the later range first so the earlier range's position stays valid).\nafterBookmark.insertText(\"\", Word.InsertLocation.replace);\nbeforeBookmark.insertText(\"\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The second paragraph of the document (the byline with the author's\n# name, course code and ID number, separated by tabs) is center-aligned\n# and stripped of all of its run content (\"Sami Mansoor Alavi\", the\n# tabs, \"BESE 8B\", the tabs, and \"209433\"), while the existing\n# \"_GoBack\" bookmark that sits in the middle of that text is preserved\n# in place.\n\n$d = $word.ActiveDocument\n\n# Paragraph 1 is the \"Secure Wallet\" title; paragraph 2 is the byline\n# paragraph that needs to be cleared and centered.\n$byline = $d.Paragraphs(2)\n\n# Center the paragraph (adds <w:jc w:val=\"center\"/> to its pPr).\n# 1 = wdAlignParagraphCenter\n$byline.Alignment = 1\n\n$bylineRange = $byline.Range\n\n# Locate the (hidden) \"_GoBack\" bookmark that lives between the tabs\n# after the name and before \"BESE 8B\" so we can delete the text around\n# it without disturbing it.\n$bookmark = $d.Bookmarks(\"_GoBack\")\n$bookmarkStart = $bookmark.Start\n$bookmarkEnd = $bookmark.End\n\n$afterBookmark = $d.Range($bookmarkEnd, $bylineRange.End)\n$beforeBookmark = $d.Range($bylineRange.Start, $bookmarkStart)\n\n# Clear the text on both sides of the bookmark (order matters: clear\n# the later range first so the earlier range's position stays valid).\n$afterBookmark.Text = \"\"\n$beforeBookmark.Text = \"\"\n"}
